$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark. In the original document it
#    sits at the end of the "Fin-Ai Web app screenshot" paragraph; after the
#    edit it needs to live at the end of the "Register page (banker..." text
#    instead, so pull it out now and re-add it in the right spot below.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) "Fin-Ai Web app" + " screenshot" -> single run "Fin-Ai Web app screenshot"
#    A whole-paragraph Find/Replace collapses the two runs into one without
#    disturbing the paragraph mark / pPr / rsid attributes.
# ---------------------------------------------------------------------------
$introRange = $d.Content
$introRange.Find.Execute(
    "Fin-Ai Web app screenshot", $true, $false, $false, $false, $false,
    $true, 1, $false, "Fin-Ai Web app screenshot", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "Register page (banker[own URL link]): adminReg.php" paragraph:
#      - merge the "Register page (", "banker[", "own URL link]): " runs
#        (and the gramStart/gramEnd proofErr markers between them) into one
#      - split "adminReg.php" into "admin" / "-" / "r" / "eg.php", wrapped in
#        a spellStart/spellEnd proofErr pair
#      - re-insert the "_GoBack" bookmark right after "eg.php" (before the
#        closing spellEnd marker)
# ---------------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Register page (banker*") {
        $targetPara = $p
        break
    }
}

$paraRange = $targetPara.Range.Duplicate

$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p w:rsidR="00C40EC9" w:rsidRDefault="00C40EC9" w:rsidP="00C40EC9">' + `
    '<w:r><w:t>Register page (banker[own URL link]): admin</w:t></w:r>' + `
    '<w:r><w:t>-</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>r</w:t></w:r>' + `
    '<w:r><w:t>eg.php</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '</w:p></w:body></w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

$paraRange.InsertXML($newParaXml)

# InsertXML drops the paragraph's pPr whenever the supplied <w:p> carries its
# own attributes, so re-apply the Heading1 style explicitly afterwards.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Register page (banker*") {
        $targetPara = $p
        break
    }
}
$targetPara.Style = "Heading 1"
